$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1444.75
$ws.Range("I2").Value = 246.2
$ws.Range("K2").Value = 246.2
$ws.Range("M2").Value = -133.2

$ws.Range("H33").Value = 7696.643
$ws.Range("I33").Value = 11327.889
$ws.Range("K33").Value = 11327.889
$ws.Range("M33").Value = -11098.889

$ws.Range("H112").Value = 1680.7778
$ws.Range("I112").Value = 1400
$ws.Range("J112").Value = 1715.875
$ws.Range("K112").Value = 4200
$ws.Range("L112").Value = 5147.625
$ws.Range("M112").Value = -3092
$ws.Range("N112").Value = -7363.625

$ws.Range("H116").Value = 11565.346
$ws.Range("J116").Value = 13942.947
$ws.Range("L116").Value = 13942.947
$ws.Range("N116").Value = -20826.947

$ws.Range("H129").Value = 1120.258
$ws.Range("I129").Value = 889
$ws.Range("K129").Value = 2667
$ws.Range("M129").Value = 2333

$ws.Range("H131").Value = 8075.5
$ws.Range("I131").Value = 921.2
$ws.Range("K131").Value = 2763.6
$ws.Range("M131").Value = 2276.4

$ws.Range("H137").Value = 1674.96
$ws.Range("I137").Value = 1205
$ws.Range("J137").Value = 2108.7693
$ws.Range("K137").Value = 3615
$ws.Range("L137").Value = 6326.3079
$ws.Range("M137").Value = -1065
$ws.Range("N137").Value = -11426.3079

$ws.Range("H138").Value = 2640.6924
$ws.Range("I138").Value = 2492.6365
$ws.Range("K138").Value = 7477.9095
$ws.Range("M138").Value = -2337.9095

$ws.Range("H141").Value = 2799.1333
$ws.Range("I141").Value = 2629.923
$ws.Range("K141").Value = 7889.768999999999
$ws.Range("M141").Value = -2709.768999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 1567.3334
$ws.Range("I19").Value = 1101
$ws.Range("J19").Value = 2500
$ws.Range("K19").Value = 1101
$ws.Range("L19").Value = 2500
$ws.Range("M19").Value = -872
$ws.Range("N19").Value = -2958

$ws.Range("H61").Value = 1940.6923
$ws.Range("I61").Value = 2039
$ws.Range("K61").Value = 2039
$ws.Range("M61").Value = -1827

$ws.Range("H74").Value = 4244.75
$ws.Range("I74").Value = 3993
$ws.Range("K74").Value = 3993
$ws.Range("M74").Value = -3119

$ws.Range("H77").Value = 4244.75
$ws.Range("I77").Value = 3993
$ws.Range("K77").Value = 19965
$ws.Range("M77").Value = -15597

$ws.Range("H102").Value = 2809.4443
$ws.Range("I102").Value = 2246.7144
$ws.Range("K102").Value = 2246.7144
$ws.Range("M102").Value = -624.7143999999998

$ws.Range("H110").Value = 2377.1
$ws.Range("I110").Value = 2085.6667
$ws.Range("K110").Value = 2085.6667
$ws.Range("M110").Value = -40.66670000000022

$ws.Range("H136").Value = 1940.6923
$ws.Range("I136").Value = 2039
$ws.Range("K136").Value = 6117
$ws.Range("M136").Value = -3567

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4377.05
$ws.Range("I86").Value = 2764.6667
$ws.Range("K86").Value = 2764.6667
$ws.Range("M86").Value = -1641.6667

$ws.Range("H89").Value = 4377.05
$ws.Range("I89").Value = 2764.6667
$ws.Range("K89").Value = 13823.3335
$ws.Range("M89").Value = -8207.333500000001

$ws.Range("H105").Value = 2991.4922
$ws.Range("I105").Value = 2794.0168
$ws.Range("K105").Value = 2794.0168
$ws.Range("M105").Value = -1047.0168

$ws.Range("H107").Value = 2601.138
$ws.Range("I107").Value = 1224.1111
$ws.Range("J107").Value = 4854.4546
$ws.Range("K107").Value = 1224.1111
$ws.Range("L107").Value = 4854.4546
$ws.Range("M107").Value = 695.8888999999999
$ws.Range("N107").Value = -8694.454600000001

$ws.Range("H119").Value = 200000
$ws.Range("J119").Value = 200000
$ws.Range("L119").Value = 200000
$ws.Range("N119").Value = -209676

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4333.1787
$ws.Range("I31").Value = 2936.5557
$ws.Range("J31").Value = 6847.1
$ws.Range("K31").Value = 2936.5557
$ws.Range("L31").Value = 6847.1
$ws.Range("M31").Value = -2641.5557
$ws.Range("N31").Value = -7437.1

$ws.Range("H34").Value = 4333.1787
$ws.Range("I34").Value = 2936.5557
$ws.Range("J34").Value = 6847.1
$ws.Range("K34").Value = 2936.5557
$ws.Range("L34").Value = 6847.1
$ws.Range("M34").Value = -2734.5557
$ws.Range("N34").Value = -7251.1

$ws.Range("H86").Value = 9811.833000000001
$ws.Range("J86").Value = 9792.5
$ws.Range("L86").Value = 9792.5
$ws.Range("N86").Value = -12038.5

$ws.Range("H89").Value = 9811.833000000001
$ws.Range("J89").Value = 9792.5
$ws.Range("L89").Value = 48962.5
$ws.Range("N89").Value = -60194.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 916.3333
$ws.Range("I5").Value = 374.5
$ws.Range("K5").Value = 1123.5
$ws.Range("M5").Value = -1011.5

$ws.Range("H40").Value = 207.77777
$ws.Range("J40").Value = 225
$ws.Range("L40").Value = 900
$ws.Range("N40").Value = -1038

$ws.Range("H129").Value = 6810.909
$ws.Range("I129").Value = 9688.083000000001
$ws.Range("K129").Value = 29064.249
$ws.Range("M129").Value = -24064.249

$ws.Range("H135").Value = 916.3333
$ws.Range("I135").Value = 374.5
$ws.Range("K135").Value = 3370.5
$ws.Range("M135").Value = -835.5

$ws.Range("H140").Value = 2441.3914
$ws.Range("I140").Value = 2150.0952
$ws.Range("K140").Value = 6450.285600000001
$ws.Range("M140").Value = -1270.285600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3198.8
$ws.Range("I102").Value = 2911
$ws.Range("J102").Value = 4350
$ws.Range("K102").Value = 2911
$ws.Range("L102").Value = 4350
$ws.Range("M102").Value = -1289
$ws.Range("N102").Value = -7594

$ws.Range("H107").Value = 67364.87
$ws.Range("I107").Value = 111263.445
$ws.Range("K107").Value = 111263.445
$ws.Range("M107").Value = -109343.445

$ws.Range("H122").Value = 3032.7
$ws.Range("I122").Value = 2466.6155
$ws.Range("K122").Value = 7399.8465
$ws.Range("M122").Value = -4949.8465

$ws.Range("H132").Value = 252249.5
$ws.Range("I132").Value = 252249.5
$ws.Range("K132").Value = 756748.5
$ws.Range("M132").Value = -754218.5

$ws.Range("H136").Value = 34195.75
$ws.Range("J136").Value = 34195.75
$ws.Range("L136").Value = 102587.25
$ws.Range("N136").Value = -107687.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").Value = ""

$ws.Range("H7").Value = 2228
$ws.Range("I7").Value = 2228
$ws.Range("K7").Value = 2228
$ws.Range("M7").Value = -2116

$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").Value = ""

$ws.Range("H100").Value = 3510.8667
$ws.Range("I100").Value = 3110.7778
$ws.Range("J100").Value = 4111
$ws.Range("K100").Value = 3110.7778
$ws.Range("L100").Value = 4111
$ws.Range("M100").Value = -2569.7778
$ws.Range("N100").Value = -5193

$ws.Range("H126").Value = 2228
$ws.Range("I126").Value = 2228
$ws.Range("K126").Value = 6684
$ws.Range("M126").Value = -4214

$ws.Range("H136").Value = 7750
$ws.Range("I136").Value = 5000
$ws.Range("J136").Value = 10500
$ws.Range("K136").Value = 15000
$ws.Range("L136").Value = 31500
$ws.Range("M136").Value = -12450
$ws.Range("N136").Value = -36600

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = ""

$ws.Range("H34").Value = 49999
$ws.Range("I34").Value = 49999
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 49999
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -49796
$ws.Range("N34").Value = ""

$ws.Range("H43").Value = 53332.668
$ws.Range("I43").Value = 50000
$ws.Range("K43").Value = 50000
$ws.Range("M43").Value = -49851

$ws.Range("H122").Value = 2425.8
$ws.Range("I122").Value = 2290.75
$ws.Range("K122").Value = 6872.25
$ws.Range("M122").Value = -4422.25

$ws.Range("H126").Value = 50597.047
$ws.Range("I126").Value = 61219.293
$ws.Range("K126").Value = 183657.879
$ws.Range("M126").Value = -181187.879

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = ""

$ws.Range("H137").Value = 95000
$ws.Range("J137").Value = 95000
$ws.Range("L137").Value = 95000
$ws.Range("N137").Value = -105200
